$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the brand_name list shared string (typos corrected, referenced by D5, D7, D8)
$fixedBrands = "ABSOLUT,JAMESON,MALIBU,GLENLIVET,KAHLUA,SEAGRAM'S,CHIVAS REGAL,AVION"
$ws.Range("D5").Value = $fixedBrands
$ws.Range("D7").Value = $fixedBrands
$ws.Range("D8").Value = $fixedBrands

# Adjust row heights
$ws.Rows.Item(5).RowHeight = 23.85
$ws.Rows.Item(7).RowHeight = 23.95
$ws.Rows.Item(8).RowHeight = 57.7

# Move the active selection to D5
$ws.Range("D5").Select()
